$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It is inserted as
# row 35 (pushing the existing rows 35-120 down to 36-121), mirroring how
# the source system prepends the latest observation to this sheet.
$ws.Rows("35").Insert()

$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value = 44838
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = 100112001
$ws.Range("G35").Value = "Berenjena"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 12000
$ws.Range("L35").Value = 13000
$ws.Range("M35").Value = 12500
$ws.Range("N35").Value = "$/caja 60 unidades"
$ws.Range("O35").Value = "Región Metropolitana"
$ws.Range("P35").Value = 208
$ws.Range("Q35").Value = 60
$ws.Range("R35").Value = "Hortaliza"
